# ---------------------------------------------------------------------------
# This script reproduces a "proofing cleanup" edit: several runs of text are
# split into multiple runs with w:proofErr spell/grammar-check bookmarks
# inserted around particular words, and a new block of text (two blank
# paragraphs followed by a paragraph of analysis text) is appended at the
# end of the document body.
#
# Because Range.Text is plain-text only, runs + w:proofErr markers have to be
# injected as raw WordOpenXML via Range.InsertXML (Flat-OPC "pkg:package"
# wrapper). InsertXML *inserts* at a collapsed range rather than replacing a
# non-collapsed one, so every paragraph we rewrite is first collapsed with
# Range.Text = "" before the replacement XML is inserted.
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

$W_NS  = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$W14_NS = 'xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"'

function Wrap-Pkg([string]$bodyXml) {
    return '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData><w:document ' + $W_NS + ' ' + $W14_NS + '>' +
        '<w:body>' + $bodyXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

# Returns the opening "<w:p ...>" tag together with its "<w:pPr>...</w:pPr>"
# (if any) for the given Paragraph, read straight off the live document so
# every rsid / paraId / textId attribute already on it is preserved exactly.
function Get-ParaHeader($para) {
    $openXml = $para.Range.WordOpenXML
    $bodyIdx = $openXml.IndexOf("<w:body>")
    $rest = $openXml.Substring($bodyIdx + 8)
    if ($rest -match '(?s)^(<w:p\b[^>]*>)(<w:pPr>.*?</w:pPr>)?') {
        return $matches[1] + $matches[2]
    }
    throw "Could not read paragraph header"
}

# Replaces the run content of $para (keeping the paragraph's own <w:p>/<w:pPr>
# untouched) with $innerXml (a sequence of <w:r>/<w:proofErr> elements).
function Replace-ParagraphInnerXml($doc, $para, [string]$innerXml) {
    $header = Get-ParaHeader $para
    $full = $para.Range
    # Leave the final paragraph-mark character alone; only collapse the part
    # that actually holds run content. This works the same whether or not
    # the paragraph is the very last one in the body.
    $sub = $doc.Range($full.Start, $full.End - 1)
    $sub.Text = ""
    $bodyXml = $header + $innerXml + "</w:p>"
    $sub.InsertXML((Wrap-Pkg $bodyXml))
}

function Find-ParagraphContaining($doc, [string]$text) {
    foreach ($p in $doc.Paragraphs) {
        if ($p.Range.Text.Contains($text)) {
            return $p
        }
    }
    throw "Paragraph containing '$text' not found"
}

$RPR_EN = '<w:rPr><w:lang w:val="en-GB"/></w:rPr>'

# ---------------------------------------------------------------------------
# 1) "Choose hypervariables " -> "Choose " + [hypervariables] + " "
# ---------------------------------------------------------------------------
$p1 = Find-ParagraphContaining $d "Choose hypervariables"
$inner1 = '<w:r w:rsidRPr="00EE6CD8">' + $RPR_EN + '<w:t xml:space="preserve">Choose </w:t></w:r>' +
          '<w:proofErr w:type="spellStart"/>' +
          '<w:r>' + $RPR_EN + '<w:t>hypervariables</w:t></w:r>' +
          '<w:proofErr w:type="spellEnd"/>' +
          '<w:r>' + $RPR_EN + '<w:t xml:space="preserve"> </w:t></w:r>'
Replace-ParagraphInnerXml $d $p1 $inner1

# ---------------------------------------------------------------------------
# 2) "Estimate Enew and performance" -> "Estimate " + [Enew] + " and performance"
# ---------------------------------------------------------------------------
$p2 = Find-ParagraphContaining $d "Estimate Enew and performance"
$inner2 = '<w:r w:rsidRPr="00EE6CD8">' + $RPR_EN + '<w:t xml:space="preserve">Estimate </w:t></w:r>' +
          '<w:proofErr w:type="spellStart"/>' +
          '<w:r>' + $RPR_EN + '<w:t>Enew</w:t></w:r>' +
          '<w:proofErr w:type="spellEnd"/>' +
          '<w:r>' + $RPR_EN + '<w:t xml:space="preserve"> and performance</w:t></w:r>'
Replace-ParagraphInnerXml $d $p2 $inner2

# ---------------------------------------------------------------------------
# 3) " Has gender balance in speaking roles changed over time (i.e. years)?"
#    -> "...time (" + [i.e.] + " years)?"
#    (the bullet-point run before the text is untouched and must be kept)
# ---------------------------------------------------------------------------
$p3 = Find-ParagraphContaining $d "Has gender balance in speaking roles changed over time"
$RPR_BULLET = '<w:rPr><w:rFonts w:ascii="NimbusRomNo9L-Regu" w:eastAsia="NimbusRomNo9L-Regu" w:cs="NimbusRomNo9L-Regu" w:hint="eastAsia"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:val="en-GB"/></w:rPr>'
$RPR_20 = '<w:rPr><w:rFonts w:ascii="NimbusRomNo9L-Regu" w:eastAsia="NimbusRomNo9L-Regu" w:cs="NimbusRomNo9L-Regu"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:val="en-GB"/></w:rPr>'
$inner3 = '<w:r w:rsidRPr="00056C10">' + $RPR_BULLET + '<w:t>' + [char]0x2022 + '</w:t></w:r>' +
          '<w:r w:rsidRPr="00056C10">' + $RPR_20 + '<w:t xml:space="preserve"> Has gender balance in speaking roles changed over time (</w:t></w:r>' +
          '<w:proofErr w:type="gramStart"/>' +
          '<w:r>' + $RPR_20 + '<w:t>i.e.</w:t></w:r>' +
          '<w:proofErr w:type="gramEnd"/>' +
          '<w:r>' + $RPR_20 + '<w:t xml:space="preserve"> years)?</w:t></w:r>'
Replace-ParagraphInnerXml $d $p3 $inner3

# ---------------------------------------------------------------------------
# 4) "I will split away maybe 15% of the date to estimate Enew. This is
#    because I can't use cross validation for this if I use cross validation
#    for feature selection..." -> split out [Enew], then append two empty
#    "NoSpacing" paragraphs and one new paragraph of analysis text at the
#    end of the document body.
# ---------------------------------------------------------------------------
$p4 = Find-ParagraphContaining $d "I will split away maybe 15% of the date"
$apos = [char]0x2019
$ellipsis = [char]0x2026
$inner4 = '<w:r>' + $RPR_EN + '<w:t xml:space="preserve">I will split away maybe 15% of the date to estimate </w:t></w:r>' +
          '<w:proofErr w:type="spellStart"/>' +
          '<w:r>' + $RPR_EN + '<w:t>Enew</w:t></w:r>' +
          '<w:proofErr w:type="spellEnd"/>' +
          '<w:r>' + $RPR_EN + '<w:t>. This is because I can' + $apos + 't use cross validation for this if I use cross validation for feature selection' + $ellipsis + '</w:t></w:r>'
Replace-ParagraphInnerXml $d $p4 $inner4

# Append the two blank paragraphs + the new analysis paragraph right after
# the (just-edited) last paragraph of the document.
$quot = [char]0x201c
$quotClose = [char]0x201d

$newBlockRuns =
    '<w:r>' + $RPR_EN + '<w:t xml:space="preserve">The analysis of the effectiveness to use </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r>' + $RPR_EN + '<w:t>kNN</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r>' + $RPR_EN + '<w:t xml:space="preserve"> to predict category ' + $quot + 'Male' + $quotClose + '</w:t></w:r>' +
    '<w:r>' + $RPR_EN + '<w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:r>' + $RPR_EN + '<w:t>or ' + $quot + 'Female in the</w:t></w:r>' +
    '<w:r>' + $RPR_EN + '<w:t xml:space="preserve"> given</w:t></w:r>' +
    '<w:r>' + $RPR_EN + '<w:t xml:space="preserve"> dataset</w:t></w:r>' +
    '<w:r>' + $RPR_EN + '<w:t xml:space="preserve"> was done in three steps. The first step included analysis of the structure of the data, in the second step the hypervariable ' + $quot + 'k' + $quotClose + ' was chosen and estimation of E[</w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r>' + $RPR_EN + '<w:t>error_new</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r>' + $RPR_EN + '<w:t xml:space="preserve">] using cross validation, in the third and final step further evaluation terms were found, also using cross validation. The result of the analysis can be seen in TABLE REF and PICTURE </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r>' + $RPR_EN + '<w:t>REF.</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>'

$newParasXml =
    '<w:p><w:pPr><w:pStyle w:val="NoSpacing"/>' + $RPR_EN + '</w:pPr></w:p>' +
    '<w:p><w:pPr><w:pStyle w:val="NoSpacing"/>' + $RPR_EN + '</w:pPr></w:p>' +
    '<w:p><w:pPr><w:pStyle w:val="NoSpacing"/>' + $RPR_EN + '</w:pPr>' + $newBlockRuns + '</w:p>'

$lastParaEnd = $p4.Range.End
$insertPoint = $d.Range($lastParaEnd - 1, $lastParaEnd - 1)
$insertPoint.InsertXML((Wrap-Pkg $newParasXml))
